$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds numeric-looking price strings (e.g. "96.497.58", "239.18") that must
# remain plain text, matching the workbook's original inlineStr cells. Forcing the
# NumberFormat to Text ("@") before assignment prevents Excel from auto-converting
# them to floating point numbers, then the original style is restored so no formatting
# side effects are introduced.

$cell = $ws.Range("D2")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "96.497.58"
$cell.Style = $origStyle
$ws.Range("E2").Value = "  +5.23%  "
$cell = $ws.Range("D3")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "3.592.67"
$cell.Style = $origStyle
$ws.Range("E3").Value = "  +8.96%  "
$ws.Range("E4").Value = "  +0.08%  "
$cell = $ws.Range("D5")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "239.18"
$cell.Style = $origStyle
$ws.Range("E5").Value = "  +4.54%  "
$cell = $ws.Range("D6")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "638.35"
$cell.Style = $origStyle
$ws.Range("E6").Value = "  +4.57%  "
$ws.Range("E7").Value = "  +7.58%  "
$cell = $ws.Range("D8")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.403"
$cell.Style = $origStyle
$ws.Range("E8").Value = "  +6.68%  "
$ws.Range("E9").Value = "  -0.07%  "
$ws.Range("E10").Value = "  +8.15%  "
$cell = $ws.Range("D11")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "3.590.56"
$cell.Style = $origStyle
$ws.Range("E11").Value = "  +9.01%  "
$ws.Range("E12").Value = "  +4.18%  "
$ws.Range("E13").Value = "  +4.65%  "
$cell = $ws.Range("D14")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "6.40"
$cell.Style = $origStyle
$ws.Range("E14").Value = "  +7.93%  "
$cell = $ws.Range("D15")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "4.266.81"
$cell.Style = $origStyle
$ws.Range("E15").Value = "  +9.27%  "
$cell = $ws.Range("D16")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "96.311.92"
$cell.Style = $origStyle
$ws.Range("E16").Value = "  +5.08%  "
$ws.Range("E17").Value = "  +5.58%  "
$cell = $ws.Range("D18")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "8.01"
$cell.Style = $origStyle
$ws.Range("E18").Value = "  +0.19%  "
$cell = $ws.Range("D19")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "3.603.41"
$cell.Style = $origStyle
$ws.Range("E19").Value = "  +9.60%  "
$cell = $ws.Range("D20")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "13.29"
$cell.Style = $origStyle
$ws.Range("E20").Value = "  +25.09%  "
$cell = $ws.Range("D21")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "18.07"
$cell.Style = $origStyle
$ws.Range("E21").Value = "  +5.97%  "
$cell = $ws.Range("D22")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.497"
$cell.Style = $origStyle
$ws.Range("E22").Value = "  +13.68%  "
$cell = $ws.Range("D23")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "516.94"
$cell.Style = $origStyle
$ws.Range("E23").Value = "  +6.70%  "
$cell = $ws.Range("D24")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "3.46"
$cell.Style = $origStyle
$ws.Range("E24").Value = "  +2.53%  "
$ws.Range("E25").Value = "  +12.17%  "
$cell = $ws.Range("D26")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "6.70"
$cell.Style = $origStyle
$ws.Range("E26").Value = "  +9.69%  "
$cell = $ws.Range("D27")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "97.22"
$cell.Style = $origStyle
$ws.Range("E27").Value = "  +9.70%  "
$cell = $ws.Range("D28")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "12.46"
$cell.Style = $origStyle
$ws.Range("E28").Value = "  +7.11%  "
$cell = $ws.Range("D29")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "3.09"
$cell.Style = $origStyle
$ws.Range("E29").Value = "  +19.32%  "
$cell = $ws.Range("D30")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "11.56"
$cell.Style = $origStyle
$ws.Range("E30").Value = "  +6.72%  "
$cell = $ws.Range("D31")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.143"
$cell.Style = $origStyle
$ws.Range("E31").Value = "  +5.21%  "
$ws.Range("E32").Value = "  -0.05%  "
$ws.Range("E33").Value = "  +7.02%  "
$cell = $ws.Range("D34")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.998"
$cell.Style = $origStyle
$ws.Range("E34").Value = "  -0.94%  "
$cell = $ws.Range("D35")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "30.42"
$cell.Style = $origStyle
$ws.Range("E35").Value = "  +9.87%  "
$ws.Range("E36").Value = "  +10.63%  "
$cell = $ws.Range("D37")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "574.09"
$cell.Style = $origStyle
$ws.Range("E37").Value = "  +7.50%  "
$cell = $ws.Range("D38")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "7.89"
$cell.Style = $origStyle
$ws.Range("E38").Value = "  +8.91%  "
$ws.Range("E39").Value = "  +10.84%  "
$ws.Range("E40").Value = "  +4.28%  "
$ws.Range("E41").Value = "  -0.01%  "
$cell = $ws.Range("D42")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.924"
$cell.Style = $origStyle
$ws.Range("E42").Value = "  +8.70%  "
$ws.Range("E43").Value = "  +6.55%  "
$ws.Range("E44").Value = "  +6.70%  "
$cell = $ws.Range("D45")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "23.79"
$cell.Style = $origStyle
$ws.Range("E45").Value = "  +0.31%  "
$cell = $ws.Range("D46")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "5.66"
$cell.Style = $origStyle
$ws.Range("E46").Value = "  +7.61%  "
$ws.Range("E47").Value = "  +7.01%  "
$ws.Range("E48").Value = "  -1.46%  "
$cell = $ws.Range("D49")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "53.82"
$cell.Style = $origStyle
$ws.Range("E49").Value = "  +4.71%  "
$ws.Range("E50").Value = "  +3.07%  "
$cell = $ws.Range("D51")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "3.12"
$cell.Style = $origStyle
$ws.Range("E51").Value = "  +5.64%  "
